# Append 9 new "reg_center_user" test-data rows (regcntr_id 10002-10010,
# usr_id 110021-110029) below the existing 21 data rows, matching the
# lang_code/is_active/cr_by/cr_dtimes pattern already used throughout the
# sheet, then mirror the trailing housekeeping Excel performs when a user
# saves after such an edit (page orientation + the "click past the last
# row" selection left on the sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# regcntr_id, usr_id pairs for the new rows (22-30)
$newRows = @(
    @(10002, 110021),
    @(10003, 110022),
    @(10004, 110023),
    @(10005, 110024),
    @(10006, 110025),
    @(10007, 110026),
    @(10008, 110027),
    @(10009, 110028),
    @(10010, 110029)
)

$startRow = 22
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $regcntrId = $newRows[$i][0]
    $usrId = $newRows[$i][1]

    $ws.Cells.Item($r, 1).Value = $regcntrId   # regcntr_id
    $ws.Cells.Item($r, 2).Value = $usrId       # usr_id
    $ws.Cells.Item($r, 3).Value = "eng"        # lang_code
    $ws.Cells.Item($r, 4).Value = $true        # is_active
    $ws.Cells.Item($r, 5).Value = "superadmin" # cr_by
    $ws.Cells.Item($r, 6).Value = "now()"      # cr_dtimes
}

# Page orientation explicitly set to portrait (adds <pageSetup orientation="portrait"/>)
$ws.PageSetup.Orientation = 1

# Leave the selection where Excel would after adding rows through row 30 and
# then selecting the next full row downward (A31:XFD1048576).
$ws.Range("A31:XFD1048576").Select()
